# Add the season "record" columns (Wins / Losses / Ties) to the player
# stats table: three new header cells in row 1 (AD1:AF1) formatted like
# the existing headers, and a value in every data row (2-43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cells the same formatting as the other headers
# (bold font, border, centered alignment) by copying an existing header
# cell's formatting, then overwrite the copied text with the real label.
$ws.Range("AA1").Copy($ws.Range("AD1"))
$ws.Range("AA1").Copy($ws.Range("AE1"))
$ws.Range("AA1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (77 wins, 85 losses, 0 ties) for every
# player row in the table.
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 30).Value = 77
    $ws.Cells.Item($row, 31).Value = 85
    $ws.Cells.Item($row, 32).Value = 0
}
